$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.566583633422852
$ws.Range("B1").Value = 1.671361207962036
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.63399076461792
$ws.Range("E1").Value = 0.7264934778213501
